# The commit swaps the contents of ppt/theme/theme1.xml ("Integral" theme,
# used by the slide master / whole deck) and ppt/theme/theme2.xml ("Office
# Theme", used only by the notes master). Apart from the <a:theme name="..">
# / <a:clrScheme name=".."> labels, the font scheme and format scheme of the
# two themes are byte-for-byte identical - the only real payload that moves
# is the 12-slot colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# So: recolour the (single, shared) editable theme so it carries the
# "Office Theme" palette that used to live in theme2.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# index -> (was "Integral" value) -> now "Office" value
$officeColors = @(
    0,            # 1  dk1      000000 (unchanged)
    16777215,     # 2  lt1      FFFFFF (unchanged)
    6968388,      # 3  dk2      44546A
    15132391,     # 4  lt2      E7E6E6
    13998939,     # 5  accent1  5B9BD5
    3243501,      # 6  accent2  ED7D31
    10855845,     # 7  accent3  A5A5A5
    49407,        # 8  accent4  FFC000
    12874308,     # 9  accent5  4472C4
    4697456,      # 10 accent6  70AD47
    12673797,     # 11 hlink    0563C1
    7491477       # 12 folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
